$d = $word.ActiveDocument

# 1. "Defendant appeared in Court on July 06, 2022" -> "... July 09, 2022"
$d.Content.Find.Execute("Defendant appeared in Court on July 06, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "Defendant appeared in Court on July 09, 2022", 2)

# 2. " license is suspended from July 06, 2022" -> " license is suspended from July 09, 2022"
$d.Content.Find.Execute(" license is suspended from July 06, 2022", $true, $false, $false, $false, $false, $true, 1, $false, " license is suspended from July 09, 2022", 2)

# 3. "September 04, 2022" -> "September 07, 2022"
$d.Content.Find.Execute("September 04, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "September 07, 2022", 2)

# 4. standalone bold "July 06, 2022" -> "July 09, 2022" (only remaining occurrence after the above replacements)
$d.Content.Find.Execute("July 06, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "July 09, 2022", 2)
